$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ticker values for columns B (Buying Opportunity) and C (support Zone) per row.
# Rows 2-3 also change the E (Short buildup) / F (FII ENTERING) tickers; rows 4-7
# have their old F ticker cleared out entirely. Rows 24-26 are brand new rows.
$data = @{
  2  = @{ B = "NSE:3PLAND";      C = "NSE:AARTIPHARM";  E = "NSE:POLYCAB"; F = "NSE:LALPATHLAB" }
  3  = @{ B = "NSE:BLUEDART";    C = "NSE:ALICON";      F = "NSE:MGL" }
  4  = @{ B = "NSE:BSE";         C = "NSE:APCOTEXIND" }
  5  = @{ B = "NSE:CEATLTD";     C = "NSE:ASHOKAMET" }
  6  = @{ B = "NSE:DPABHUSHAN";  C = "NSE:BALAMINES" }
  7  = @{ B = "NSE:GREENPANEL";  C = "NSE:BALAXI" }
  8  = @{ B = "NSE:ICICIB22";    C = "NSE:BBTCL" }
  9  = @{ B = "NSE:IGL";         C = "NSE:CINELINE" }
  10 = @{ B = "NSE:INDTERRAIN";  C = "NSE:COMPUSOFT" }
  11 = @{ B = "NSE:LALPATHLAB";  C = "NSE:ELECTCAST" }
  12 = @{ B = "NSE:MGL";         C = "NSE:FCSSOFT" }
  13 = @{ B = "NSE:MONIFTY500";  C = "NSE:GREENPOWER" }
  14 = @{ B = "NSE:RELIANCE";    C = "NSE:GUJRAFFIA" }
  15 = @{ B = "NSE:RVNL";        C = "NSE:HEMIPROP" }
  16 = @{ B = "";                C = "NSE:HERITGFOOD" }
  17 = @{ B = "";                C = "NSE:JAGRAN" }
  18 = @{ B = "";                C = "NSE:JAYNECOIND" }
  19 = @{ B = "";                C = "NSE:JBMA" }
  20 = @{ B = "";                C = "NSE:KOTHARIPRO" }
  21 = @{ B = "";                C = "NSE:LEXUS" }
  22 = @{ B = "";                C = "NSE:LXCHEM" }
  23 = @{ B = "";                C = "NSE:MASFIN" }
  24 = @{ B = "";                C = "NSE:MAYURUNIQ" }
  25 = @{ B = "";                C = "NSE:ORIENTLTD" }
  26 = @{ B = "";                C = "NSE:PNCINFRA" }
}

# First, stamp the brand-new rows 24-26 with the same formatting used by the
# existing numbered rows, by copying row 23's cells down (this brings the
# "A" column bold/centered/bordered style, and blank D/E/F cells, along).
$ws.Range("A23").Copy($ws.Range("A24:A26")) | Out-Null
$ws.Range("D23:F23").Copy($ws.Range("D24:F26")) | Out-Null

# Clear out the old F column values for rows 4-7 (no longer populated).
foreach ($r in 4..7) {
  $ws.Range("F$r").Value = ""
}

foreach ($r in 2..26) {
  $row = $data[$r]

  $ws.Range("A$r").Value = $r - 2
  if ($row.ContainsKey("B")) {
    $ws.Range("B$r").Value = $row.B
  }
  if ($row.ContainsKey("C")) {
    $ws.Range("C$r").Value = $row.C
  }
  if ($row.ContainsKey("E")) {
    $ws.Range("E$r").Value = $row.E
  }
  if ($row.ContainsKey("F")) {
    $ws.Range("F$r").Value = $row.F
  }
}
